$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.072.40'
$ws.Range("E2").Value = '  -2.20%  '
$ws.Range("D3").Value = '2.661.76'
$ws.Range("E3").Value = '  -0.88%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = "'525.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.77%  '
$ws.Range("D6").Value = "'144.47"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.06%  '
$ws.Range("E7").Value = '  +0.22%  '
$ws.Range("E8").Value = '  -0.91%  '
$ws.Range("E9").Value = '  +7.95%  '
$ws.Range("E10").Value = '  -2.23%  '
$ws.Range("E11").Value = '  -1.78%  '
$ws.Range("E12").Value = '  +1.37%  '
$ws.Range("D13").Value = '3.133.04'
$ws.Range("E13").Value = '  -0.91%  '
$ws.Range("D14").Value = '59.058.93'
$ws.Range("E14").Value = '  -2.28%  '
$ws.Range("D15").Value = "'21.11"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.60%  '
$ws.Range("D16").Value = '2.669.85'
$ws.Range("E16").Value = '  -1.18%  '
$ws.Range("E17").Value = '  -1.77%  '
$ws.Range("D18").Value = "'338.72"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.65%  '
$ws.Range("E19").Value = '  -3.41%  '
$ws.Range("D20").Value = "'10.40"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.46%  '
$ws.Range("D21").Value = "'6.39"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.14%  '
$ws.Range("D22").Value = "'0.996"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.30%  '
$ws.Range("D23").Value = "'64.36"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.32%  '
$ws.Range("D24").Value = "'0.419"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.69%  '
$ws.Range("E25").Value = '  -1.76%  '
$ws.Range("D26").Value = "'0.997"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.34%  '
$ws.Range("D27").Value = '0.0₃0801'
$ws.Range("E27").Value = '  -1.46%  '
$ws.Range("D28").Value = "'7.09"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.40%  '
$ws.Range("D29").Value = "'6.68"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.90%  '
$ws.Range("E30").Value = '  +0.08%  '
$ws.Range("D31").Value = "'1.59"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.16%  '
$ws.Range("D32").Value = "'18.86"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.05%  '
$ws.Range("D33").Value = "'150.84"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.51%  '
$ws.Range("D34").Value = "'4.15"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.07%  '
$ws.Range("E35").Value = '  -3.42%  '
$ws.Range("D36").Value = "'0.890"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -6.01%  '
$ws.Range("D37").Value = "'0.873"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.37%  '
$ws.Range("D38").Value = "'36.90"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.27%  '
$ws.Range("D39").Value = "'1.45"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -5.95%  '
$ws.Range("D40").Value = "'3.58"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.99%  '
$ws.Range("E41").Value = '  +0.68%  '
$ws.Range("E42").Value = '  +0.33%  '
$ws.Range("D43").Value = "'19.92"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.57%  '
$ws.Range("D44").Value = "'275.24"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.05%  '
$ws.Range("E45").Value = '  -1.77%  '
$ws.Range("E46").Value = '  +1.98%  '
$ws.Range("D47").Value = "'0.0531"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.55%  '
$ws.Range("D48").Value = '2.045.72'
$ws.Range("E48").Value = '  -3.83%  '
$ws.Range("B49").Value = 'VeChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D49").Value = "'0.0229"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.76%  '
$ws.Range("B50").Value = 'RenderToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D50").Value = "'4.70"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.22%  '
$ws.Range("D51").Value = "'18.92"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.44%  '
